# Upgrade database schema migration
# Appends a new row (row 29) of data to each of the four worksheets,
# mirroring the existing row layout (time, lengths, ids, checksums, and
# their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rowTime = 45815.4618287037

$sheetsData = @(
    @{ Name = "MID_LFT_#1"; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x88"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 392; I = 7 },
    @{ Name = "MID_LFT_#2"; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x78"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 376; I = 25 },
    @{ Name = "MID_PLT_#1"; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x6D"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 109; I = 15 },
    @{ Name = "MID_PLT_#2"; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x81"; E = "0x9";  F = 130; G = "5.68631262647113e+23"; H = 129; I = 9 }
)

foreach ($sheetData in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetData.Name)

    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $rowTime
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $sheetData.B
    $ws.Cells.Item($newRow, 3).Value = $sheetData.C
    $ws.Cells.Item($newRow, 4).Value = $sheetData.D
    $ws.Cells.Item($newRow, 5).Value = $sheetData.E
    $ws.Cells.Item($newRow, 6).Value = $sheetData.F
    $ws.Cells.Item($newRow, 7).Value = [double]$sheetData.G
    $ws.Cells.Item($newRow, 8).Value = $sheetData.H
    $ws.Cells.Item($newRow, 9).Value = $sheetData.I
}
